$d = $word.ActiveDocument

# Locate the paragraph that reads "Edison Achalma" and uses the "Author"
# style (the byline right under the article title).
$target = $null
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Edison Achalma" -and $p.Style.NameLocal -eq "Author") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Range.End for a paragraph points right after its own paragraph mark,
    # i.e. the start of the following paragraph's content.
    $insertPos = $target.Range.End
    $d.Range($insertPos, $insertPos).InsertParagraphAfter()

    # The newly created (empty) paragraph now starts exactly at $insertPos.
    $newPara = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -eq $insertPos) {
            $newPara = $p
            break
        }
    }

    $newPara.Range.Text = "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga"
    $newPara.Style = $d.Styles.Item("Author")
}
